$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C8").Value = "No"
$ws.Range("C9").Value = "NO"

$ws.Range("A10").Value = "Tools"
$ws.Range("B10").Value = "Admin Tools Scenarios"
$ws.Range("C10").Value = "Yes"

$ws.Range("A13").Select()
